$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 707
$ws.Range("I6").Value = 707
$ws.Range("K6").Value = 2121
$ws.Range("M6").Value = -2009
$ws.Range("H13").Value = 4000
$ws.Range("J13").Value = 4000
$ws.Range("L13").Value = 4000
$ws.Range("N13").Value = -4338
$ws.Range("H28").Value = 654.3333
$ws.Range("I28").Value = 708.3333
$ws.Range("J28").Value = 438.33334
$ws.Range("K28").Value = 708.3333
$ws.Range("L28").Value = 438.33334
$ws.Range("M28").Value = -223.3333
$ws.Range("N28").Value = -1408.33334
$ws.Range("H32").Value = 31251568
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 31251568
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 31251568
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -31252220
$ws.Range("H33").Value = 9640776
$ws.Range("I33").Value = 4839093.5
$ws.Range("K33").Value = 4839093.5
$ws.Range("M33").Value = -4838864.5
$ws.Range("H40").Value = 1949.4
$ws.Range("J40").Value = 1986.75
$ws.Range("L40").Value = 1986.75
$ws.Range("N40").Value = -2336.75
$ws.Range("H138").Value = 6485.9116
$ws.Range("J138").Value = 5190.273
$ws.Range("L138").Value = 15570.819
$ws.Range("N138").Value = -25850.819
$ws.Range("H141").Value = 1807.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 188388.72
$ws.Range("I32").Value = 223852.94
$ws.Range("J32").Value = 11067.667
$ws.Range("K32").Value = 223852.94
$ws.Range("L32").Value = 11067.667
$ws.Range("M32").Value = -223565.94
$ws.Range("N32").Value = -11641.667
$ws.Range("H110").Value = 1110.7142
$ws.Range("I110").Value = 971.13635
$ws.Range("J110").Value = 1622.5
$ws.Range("K110").Value = 971.13635
$ws.Range("L110").Value = 1622.5
$ws.Range("M110").Value = 1073.86365
$ws.Range("N110").Value = -5712.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1037.5333
$ws.Range("I22").Value = 727.9231
$ws.Range("K22").Value = 727.9231
$ws.Range("M22").Value = -554.9231
$ws.Range("H86").Value = 7209.5264
$ws.Range("I86").Value = 5199.6
$ws.Range("K86").Value = 5199.6
$ws.Range("M86").Value = -4076.6
$ws.Range("H89").Value = 7209.5264
$ws.Range("I89").Value = 5199.6
$ws.Range("K89").Value = 25998
$ws.Range("M89").Value = -20382
$ws.Range("H99").Value = 9313.166999999999
$ws.Range("I99").Value = 10102.3125
$ws.Range("K99").Value = 10102.3125
$ws.Range("M99").Value = -8604.3125
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2094.0232
$ws.Range("I31").Value = 1891.2927
$ws.Range("J31").Value = 6250
$ws.Range("K31").Value = 1891.2927
$ws.Range("L31").Value = 6250
$ws.Range("M31").Value = -1596.2927
$ws.Range("N31").Value = -6840
$ws.Range("H34").Value = 2094.0232
$ws.Range("I34").Value = 1891.2927
$ws.Range("J34").Value = 6250
$ws.Range("K34").Value = 1891.2927
$ws.Range("L34").Value = 6250
$ws.Range("M34").Value = -1689.2927
$ws.Range("N34").Value = -6654
$ws.Range("H62").Value = 6375
$ws.Range("J62").Value = 7166.6665
$ws.Range("L62").Value = 7166.6665
$ws.Range("N62").Value = -8414.666499999999
$ws.Range("H65").Value = 6375
$ws.Range("J65").Value = 7166.6665
$ws.Range("L65").Value = 35833.3325
$ws.Range("N65").Value = -42073.3325
$ws.Range("H99").Value = 4001300
$ws.Range("J99").Value = 1100
$ws.Range("L99").Value = 1100
$ws.Range("N99").Value = -4096
$ws.Range("H102").Value = 50000.5
$ws.Range("J102").Value = 50000.5
$ws.Range("L102").Value = 50000.5
$ws.Range("N102").Value = -54868.5
$ws.Range("H126").Value = 4001300
$ws.Range("J126").Value = 1100
$ws.Range("L126").Value = 3300
$ws.Range("N126").Value = -8240
$ws.Range("H141").Value = 369999.8
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 369999.8
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 369999.8
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -380359.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 342.6
$ws.Range("I98").Value = 169.5
$ws.Range("J98").Value = 458
$ws.Range("K98").Value = 508.5
$ws.Range("L98").Value = 1374
$ws.Range("M98").Value = 989.5
$ws.Range("N98").Value = -4370
$ws.Range("H125").Value = 5100
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("H129").Value = 6059.4165
$ws.Range("J129").Value = 9403.571
$ws.Range("L129").Value = 28210.713
$ws.Range("N129").Value = -38210.713
$ws.Range("H131").Value = 2676814
$ws.Range("J131").Value = 3528
$ws.Range("L131").Value = 10584
$ws.Range("N131").Value = -20664
$ws.Range("H138").Value = 3272.3684
$ws.Range("I138").Value = 3367.889
$ws.Range("J138").Value = 1553
$ws.Range("K138").Value = 10103.667
$ws.Range("L138").Value = 4659
$ws.Range("M138").Value = -4963.667000000001
$ws.Range("N138").Value = -14939
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1999.5
$ws.Range("I29").Value = 1999
$ws.Range("K29").Value = 1999
$ws.Range("M29").Value = -1709
$ws.Range("H70").Value = 3381.0908
$ws.Range("I70").Value = 3314.4285
$ws.Range("K70").Value = 3314.4285
$ws.Range("M70").Value = -3044.4285
$ws.Range("H73").Value = 3381.0908
$ws.Range("I73").Value = 3314.4285
$ws.Range("K73").Value = 3314.4285
$ws.Range("M73").Value = -2378.4285
$ws.Range("H132").Value = 2566090.8
$ws.Range("I132").Value = 1694.2858
$ws.Range("K132").Value = 5082.857400000001
$ws.Range("M132").Value = -2552.857400000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2155.3462
$ws.Range("I40").Value = 2343.4375
$ws.Range("J40").Value = 1854.4
$ws.Range("K40").Value = 2343.4375
$ws.Range("L40").Value = 1854.4
$ws.Range("M40").Value = -2207.4375
$ws.Range("N40").Value = -2126.4
$ws.Range("H55").Value = 431.59375
$ws.Range("I55").Value = 319.8125
$ws.Range("J55").Value = 543.375
$ws.Range("K55").Value = 319.8125
$ws.Range("L55").Value = 543.375
$ws.Range("M55").Value = -146.8125
$ws.Range("N55").Value = -889.375
$ws.Range("H61").Value = 16809
$ws.Range("I61").Value = 18022.5
$ws.Range("K61").Value = 18022.5
$ws.Range("M61").Value = -17820.5
$ws.Range("H82").Value = 994.1667
$ws.Range("I82").Value = 1074
$ws.Range("K82").Value = 1074
$ws.Range("M82").Value = -713
$ws.Range("H85").Value = 994.1667
$ws.Range("I85").Value = 1074
$ws.Range("K85").Value = 1074
$ws.Range("M85").Value = 174
$ws.Range("H93").Value = 1045.5
$ws.Range("I93").Value = 1045.5
$ws.Range("K93").Value = 1045.5
$ws.Range("M93").Value = 202.5
$ws.Range("H113").Value = 16809
$ws.Range("I113").Value = 18022.5
$ws.Range("K113").Value = 18022.5
$ws.Range("M113").Value = -15852.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2999.5
$ws.Range("I62").Value = 2499.5
$ws.Range("J62").Value = 3499.5
$ws.Range("K62").Value = 2499.5
$ws.Range("L62").Value = 3499.5
$ws.Range("M62").Value = -1875.5
$ws.Range("N62").Value = -4747.5
$ws.Range("H65").Value = 2999.5
$ws.Range("I65").Value = 2499.5
$ws.Range("J65").Value = 3499.5
$ws.Range("K65").Value = 12497.5
$ws.Range("L65").Value = 17497.5
$ws.Range("M65").Value = -9377.5
$ws.Range("N65").Value = -23737.5
$ws.Range("H96").Value = 3010.5293
$ws.Range("J96").Value = 3165.5
$ws.Range("L96").Value = 3165.5
$ws.Range("N96").Value = -5911.5
$ws.Range("H126").Value = 3596.6562
$ws.Range("I126").Value = 3520.6365
$ws.Range("J126").Value = 3763.9
$ws.Range("K126").Value = 10561.9095
$ws.Range("L126").Value = 11291.7
$ws.Range("M126").Value = -8091.9095
$ws.Range("N126").Value = -16231.7
